$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 32, shifting existing rows (old 32..332) down to (33..333)
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new translation-key entry
$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = "In ""settings"" form, tab ""User interface"""
$ws.Range("E32").Value = "Remember window position and size on startup"
$ws.Rows("32:32").RowHeight = 30

# Existing row 25 (strChkDlgPath) gains the same Comment text and row height
$ws.Range("D25").Value = "In ""settings"" form, tab ""User interface"""
$ws.Rows("25:25").RowHeight = 30

# Grow the "Tabla13" table/autofilter range by one row to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))
